# Update the "quantity" column (F) values on the NetDemand sheet
# to reflect use of uncon_planned_qty for future production while
# keeping produced quantity for today.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NetDemand")

$ws.Range("F2").Value = -239
$ws.Range("F3").Value = -530
$ws.Range("F4").Value = -832
$ws.Range("F5").Value = -107
$ws.Range("F7").Value = -107
